$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BOM table gained a new component row: an inductor (L1, 33uH) added
# for ripple reduction on the updated PCB layout.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

$newRow.Range.Item(1, 1).Value = 1
$newRow.Range.Item(1, 2).Value = "L"
$newRow.Range.Item(1, 3).Value = "33uH"
$newRow.Range.Item(1, 4).Value = "L1"
$newRow.Range.Item(1, 5).Value = "L-PIHV4119 33µ"

# Match the author's final selection before saving.
$ws.Range("E20").Select()
